# Update "想去人数" (F column) values across sheets to reflect the latest
# scrape snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 240
$ws1.Range("F3").Value = 4893
$ws1.Range("F8").Value = 102
$ws1.Range("F19").Value = 4226
$ws1.Range("F20").Value = 6525
$ws1.Range("F22").Value = 43
$ws1.Range("F23").Value = 89
$ws1.Range("F44").Value = 510

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 115

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 240
$ws4.Range("F3").Value = 4893
$ws4.Range("F7").Value = 115
$ws4.Range("F9").Value = 102
$ws4.Range("F20").Value = 4226
$ws4.Range("F21").Value = 6525
$ws4.Range("F23").Value = 43
$ws4.Range("F24").Value = 89
$ws4.Range("F45").Value = 510
